$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 72.79331566666666
$ws.Range("H2").Value = 218.379947
$ws.Range("I2").Value = 0.2828741606141505
$ws.Range("J2").Value = 0.2828741606141506
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.502995666666665
$ws.Range("N2").Value = 25.508987
$ws.Range("O2").Value = 0.7300889802301167
$ws.Range("P2").Value = 0.7300889802301168
$ws.Range("Q2").Value = 618.9612476759653
$ws.Range("R2").Value = 5570.651229083687
$ws.Range("S2").Value = 0.2065233074562354
$ws.Range("T2").Value = 0.2065233074562355

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 72.79331566666666
$ws.Range("H3").Value = 218.379947
$ws.Range("I3").Value = 0.2828741606141505
$ws.Range("J3").Value = 0.2828741606141506
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.074135666666667
$ws.Range("N3").Value = 3.222407
$ws.Range("O3").Value = 0.09222803871107818
$ws.Range("P3").Value = 0.0922280387110782
$ws.Range("Q3").Value = 78.1898966524921
$ws.Range("R3").Value = 703.7090698724289
$ws.Range("S3").Value = 0.02608892903548563
$ws.Range("T3").Value = 0.02608892903548563

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 72.79331566666666
$ws.Range("H4").Value = 218.379947
$ws.Range("I4").Value = 0.2828741606141505
$ws.Range("J4").Value = 0.2828741606141506
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.069388333333333
$ws.Range("N4").Value = 6.208165
$ws.Range("O4").Value = 0.177682981058805
$ws.Range("P4").Value = 0.177682981058805
$ws.Range("Q4").Value = 150.6376381852506
$ws.Range("R4").Value = 1355.738743667255
$ws.Range("S4").Value = 0.05026192412242948
$ws.Range("T4").Value = 0.0502619241224295

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 117.1700846666667
$ws.Range("H5").Value = 351.510254
$ws.Range("I5").Value = 0.4553218801152877
$ws.Range("J5").Value = 0.4553218801152878
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.502995666666665
$ws.Range("N5").Value = 25.508987
$ws.Range("O5").Value = 0.7300889802301167
$ws.Range("P5").Value = 0.7300889802301168
$ws.Range("Q5").Value = 996.296722183633
$ws.Range("R5").Value = 8966.670499652699
$ws.Range("S5").Value = 0.3324254871298298
$ws.Range("T5").Value = 0.3324254871298299

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 117.1700846666667
$ws.Range("H6").Value = 351.510254
$ws.Range("I6").Value = 0.4553218801152877
$ws.Range("J6").Value = 0.4553218801152878
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.074135666666667
$ws.Range("N6").Value = 3.222407
$ws.Range("O6").Value = 0.09222803871107818
$ws.Range("P6").Value = 0.0922280387110782
$ws.Range("Q6").Value = 125.8565670068198
$ws.Range("R6").Value = 1132.709103061378
$ws.Range("S6").Value = 0.04199344398527365
$ws.Range("T6").Value = 0.04199344398527367

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 117.1700846666667
$ws.Range("H7").Value = 351.510254
$ws.Range("I7").Value = 0.4553218801152877
$ws.Range("J7").Value = 0.4553218801152878
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.069388333333333
$ws.Range("N7").Value = 6.208165
$ws.Range("O7").Value = 0.177682981058805
$ws.Range("P7").Value = 0.177682981058805
$ws.Range("Q7").Value = 242.4704062248789
$ws.Range("R7").Value = 2182.23365602391
$ws.Range("S7").Value = 0.08090294900018415
$ws.Range("T7").Value = 0.08090294900018417

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 67.37122333333333
$ws.Range("H8").Value = 202.11367
$ws.Range("I8").Value = 0.2618039592705617
$ws.Range("J8").Value = 0.2618039592705618
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 8.502995666666665
$ws.Range("N8").Value = 25.508987
$ws.Range("O8").Value = 0.7300889802301167
$ws.Range("P8").Value = 0.7300889802301168
$ws.Range("Q8").Value = 572.8572200613654
$ws.Range("R8").Value = 5155.714980552289
$ws.Range("S8").Value = 0.1911401856440514
$ws.Range("T8").Value = 0.1911401856440515

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 67.37122333333333
$ws.Range("H9").Value = 202.11367
$ws.Range("I9").Value = 0.2618039592705617
$ws.Range("J9").Value = 0.2618039592705618
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.074135666666667
$ws.Range("N9").Value = 3.222407
$ws.Range("O9").Value = 0.09222803871107818
$ws.Range("P9").Value = 0.0922280387110782
$ws.Range("Q9").Value = 72.36583388929888
$ws.Range("R9").Value = 651.29250500369
$ws.Range("S9").Value = 0.0241456656903189
$ws.Range("T9").Value = 0.02414566569031891

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 67.37122333333333
$ws.Range("H10").Value = 202.11367
$ws.Range("I10").Value = 0.2618039592705617
$ws.Range("J10").Value = 0.2618039592705618
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.069388333333333
$ws.Range("N10").Value = 6.208165
$ws.Range("O10").Value = 0.177682981058805
$ws.Range("P10").Value = 0.177682981058805
$ws.Range("Q10").Value = 139.4172235683944
$ws.Range("R10").Value = 1254.75501211555
$ws.Range("S10").Value = 0.04651810793619138
$ws.Range("T10").Value = 0.04651810793619139

